$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy formatting (styles/number formats) from the last existing data row (223)
# down into the 4 new rows, then drop the stray empty column-E cells that
# PasteSpecial creates when it spans the whole row.
$ws.Range("A223:AK223").Copy() | Out-Null
$ws.Range("A224:AK227").PasteSpecial(-4122) | Out-Null
$ws.Range("E224:E227").Clear() | Out-Null

# --- Row 224 ---
$ws.Cells.Item(224, 1).Value = 54
$ws.Cells.Item(224, 2).Value = 44877
$ws.Cells.Item(224, 3).Value = "Volkan"
$ws.Cells.Item(224, 4).Value = 1
$ws.Cells.Item(224, 6).Value = "Batuhan"
$ws.Cells.Item(224, 7).Value = 9
$ws.Cells.Item(224, 8).Value = "Mavi"
$ws.Cells.Item(224, 9).Value = 3
$ws.Cells.Item(224, 10).Value = 1
$ws.Cells.Item(224, 11).Value = 5
$ws.Cells.Item(224, 12).Value = 5
$ws.Cells.Item(224, 13).Value = 0
$ws.Cells.Item(224, 14).Value = 1
$ws.Cells.Item(224, 15).Value = 0
$ws.Cells.Item(224, 16).Value = 2
$ws.Cells.Item(224, 17).Value = "Kaya, Koyun, Koyun"
$ws.Cells.Item(224, 18).Value = 0
$ws.Cells.Item(224, 19).Value = 80
$ws.Cells.Item(224, 20).Value = 0
$ws.Cells.Item(224, 21).Value = 0.5590277777777778
$ws.Cells.Item(224, 22).Value = 0
$ws.Cells.Item(224, 23).Value = 0
$ws.Cells.Item(224, 24).Value = 0
$ws.Cells.Item(224, 25).Value = 0
$ws.Cells.Item(224, 26).Value = 0
$ws.Cells.Item(224, 27).Value = 0
$ws.Cells.Item(224, 28).Value = 1
$ws.Cells.Item(224, 29).Value = 0
$ws.Cells.Item(224, 30).Value = 0
$ws.Cells.Item(224, 31).Value = 0
$ws.Cells.Item(224, 32).Value = 0
$ws.Cells.Item(224, 33).Value = 1
$ws.Cells.Item(224, 34).Value = 1
$ws.Cells.Item(224, 35).Value = 0
$ws.Cells.Item(224, 36).Value = "5, 10, 6, 11, 4"
$ws.Cells.Item(224, 37).Value = "Orta Çember"

# --- Row 225 ---
$ws.Cells.Item(225, 1).Value = 54
$ws.Cells.Item(225, 2).Value = 44877
$ws.Cells.Item(225, 3).Value = "Ecem"
$ws.Cells.Item(225, 4).Value = 2
$ws.Cells.Item(225, 6).Value = "Batuhan"
$ws.Cells.Item(225, 7).Value = 9
$ws.Cells.Item(225, 8).Value = "Turuncu"
$ws.Cells.Item(225, 9).Value = 1
$ws.Cells.Item(225, 10).Value = 4
$ws.Cells.Item(225, 11).Value = 10
$ws.Cells.Item(225, 12).Value = 3
$ws.Cells.Item(225, 13).Value = 1
$ws.Cells.Item(225, 14).Value = 0
$ws.Cells.Item(225, 15).Value = 0
$ws.Cells.Item(225, 16).Value = 3
$ws.Cells.Item(225, 17).Value = "Odun, Tuğla"
$ws.Cells.Item(225, 18).Value = 0
$ws.Cells.Item(225, 19).Value = 80
$ws.Cells.Item(225, 20).Value = 0
$ws.Cells.Item(225, 21).Value = 0.5590277777777778
$ws.Cells.Item(225, 22).Value = 1
$ws.Cells.Item(225, 23).Value = 0
$ws.Cells.Item(225, 24).Value = 0
$ws.Cells.Item(225, 25).Value = 0
$ws.Cells.Item(225, 26).Value = 0
$ws.Cells.Item(225, 27).Value = 0
$ws.Cells.Item(225, 28).Value = 0
$ws.Cells.Item(225, 29).Value = 1
$ws.Cells.Item(225, 30).Value = 0
$ws.Cells.Item(225, 31).Value = 0
$ws.Cells.Item(225, 32).Value = 2
$ws.Cells.Item(225, 33).Value = 0
$ws.Cells.Item(225, 34).Value = 0
$ws.Cells.Item(225, 35).Value = 1
$ws.Cells.Item(225, 36).Value = "5, 6, 11, 3, 8"
$ws.Cells.Item(225, 37).Value = "Orta Çember"

# --- Row 226 ---
$ws.Cells.Item(226, 1).Value = 54
$ws.Cells.Item(226, 2).Value = 44877
$ws.Cells.Item(226, 3).Value = "Batuhan"
$ws.Cells.Item(226, 4).Value = 3
$ws.Cells.Item(226, 6).Value = "Batuhan"
$ws.Cells.Item(226, 7).Value = 10
$ws.Cells.Item(226, 8).Value = "Kırmızı"
$ws.Cells.Item(226, 9).Value = 4
$ws.Cells.Item(226, 10).Value = 1
$ws.Cells.Item(226, 11).Value = 6
$ws.Cells.Item(226, 12).Value = 2
$ws.Cells.Item(226, 13).Value = 0
$ws.Cells.Item(226, 14).Value = 0
$ws.Cells.Item(226, 15).Value = 0
$ws.Cells.Item(226, 16).Value = 1
$ws.Cells.Item(226, 17).Value = "Saman, Saman, Kaya"
$ws.Cells.Item(226, 18).Value = 0
$ws.Cells.Item(226, 19).Value = 80
$ws.Cells.Item(226, 20).Value = 0
$ws.Cells.Item(226, 21).Value = 0.5590277777777778
$ws.Cells.Item(226, 22).Value = 0
$ws.Cells.Item(226, 23).Value = 0
$ws.Cells.Item(226, 24).Value = 0
$ws.Cells.Item(226, 25).Value = 0
$ws.Cells.Item(226, 26).Value = 0
$ws.Cells.Item(226, 27).Value = 0
$ws.Cells.Item(226, 28).Value = 0
$ws.Cells.Item(226, 29).Value = 0
$ws.Cells.Item(226, 30).Value = 0
$ws.Cells.Item(226, 31).Value = 0
$ws.Cells.Item(226, 32).Value = 1
$ws.Cells.Item(226, 33).Value = 0
$ws.Cells.Item(226, 34).Value = 0
$ws.Cells.Item(226, 35).Value = 0
$ws.Cells.Item(226, 36).Value = "5, 9, 10, 4, 9, 11"
$ws.Cells.Item(226, 37).Value = "Orta Çember"

# --- Row 227 ---
$ws.Cells.Item(227, 1).Value = 54
$ws.Cells.Item(227, 2).Value = 44877
$ws.Cells.Item(227, 3).Value = "Alperen"
$ws.Cells.Item(227, 4).Value = 4
$ws.Cells.Item(227, 6).Value = "Batuhan"
$ws.Cells.Item(227, 7).Value = 9
$ws.Cells.Item(227, 8).Value = "Beyaz"
$ws.Cells.Item(227, 9).Value = 3
$ws.Cells.Item(227, 10).Value = 1
$ws.Cells.Item(227, 11).Value = 6
$ws.Cells.Item(227, 12).Value = 4
$ws.Cells.Item(227, 13).Value = 0
$ws.Cells.Item(227, 14).Value = 0
$ws.Cells.Item(227, 15).Value = 1
$ws.Cells.Item(227, 16).Value = 0
$ws.Cells.Item(227, 17).Value = "Tuğla, Odun, Kaya"
$ws.Cells.Item(227, 18).Value = 1
$ws.Cells.Item(227, 19).Value = 80
$ws.Cells.Item(227, 20).Value = 0
$ws.Cells.Item(227, 21).Value = 0.5590277777777778
$ws.Cells.Item(227, 22).Value = 0
$ws.Cells.Item(227, 23).Value = 1
$ws.Cells.Item(227, 24).Value = 0
$ws.Cells.Item(227, 25).Value = 0
$ws.Cells.Item(227, 26).Value = 0
$ws.Cells.Item(227, 27).Value = 0
$ws.Cells.Item(227, 28).Value = 0
$ws.Cells.Item(227, 29).Value = 0
$ws.Cells.Item(227, 30).Value = 0
$ws.Cells.Item(227, 31).Value = 0
$ws.Cells.Item(227, 32).Value = 0
$ws.Cells.Item(227, 33).Value = 0
$ws.Cells.Item(227, 34).Value = 0
$ws.Cells.Item(227, 35).Value = 0
$ws.Cells.Item(227, 36).Value = "3, 8, 10, 3, 4, 19"
$ws.Cells.Item(227, 37).Value = "Orta Çember"
